$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kingdoms")

# Row 4 (Ireland) was missing the "feudal_government" value in the
# "government" column (E) that the other kingdoms (rows 2-3) already have.
$ws.Range("E4").Value = "feudal_government"

# Move the saved selection/active cell to H4, matching the latest view state.
$ws.Range("H4").Select() | Out-Null
